$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(10, "aa", "dssdfd"),
    @(11, "bb", "gfdsgdfs"),
    @(12, "cc", "sfdgf"),
    @(13, "dd", "sdfgfsd"),
    @(14, "ee", "sdfgsfd"),
    @(15, "ff", "sfgsfg"),
    @(16, "gg", "fsdgdfs"),
    @(17, "hh", "sfgfsdfg"),
    @(18, "ii", "sdfg"),
    @(19, "jj", "sfgffd")
)

$row = 11
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}
